# edit.ps1
# Applies the "audit reset / no-violations" update to the weekly report workbook:
#  - Updates the "Report Generated On" timestamp
#  - Zeroes out the Total Billed Amount summary cell
#  - Zeroes out every daily-detail and daily-total "Pricing" value in column H

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report-generated timestamp shown near the top of the sheet.
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# Reset the summary "Total Billed Amount" figure.
$ws.Range("C8").Value = 0

# Zero every pricing value (column H) across all daily detail rows and totals
# affected by the audit re-run.
$rowsToZero = @(16,17,22,23,24,25,26,27,28,29,30,31,32,37,38,39,40,41,42,43,48,49,50,51,52,57,58,59,60,61,62,63,64)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 8).Value = 0
}
